$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial for every data row (2-129).
# The update bumps that date forward by one day (46060 -> 46061) for all rows.
$lastRow = $ws.UsedRange.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $current = $cell.Value()
    if ($current -ne $null) {
        $cell.Value = $current.AddDays(1)
    }
}
